$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Remove the stray duplicate chart-data defined names left over from the
# chart being recreated/refreshed (v1.0/v1.1 are still used by chartEx1.xml)
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# New header labels (bold, matching the style of the other stat headers
# such as D6/D9/D12 "Min"/"Q1"/"Q3")
$ws.Range("D15").Value = "Std"
$ws.Range("E15").Value = "Relative std"
$ws.Range("D15:E15").Font.Bold = $true

# New computed statistics
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "= (D16 / E3) * 100"

# Update the active selection to match the edited area
$ws.Range("E17").Select()
